$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old hand-rolled "Id;..." CSV-style BOM stub with a proper
# KiCad-exported Bill of Materials grid (grouped by value/footprint).
$ws.Cells.Clear()

# Row 1: component count summary
$ws.Range("A1").Value = "Component Count:"
$ws.Range("B1").Value = 8

# Row 3: column headers
$ws.Range("A3").Value = "Ref"
$ws.Range("B3").Value = "Qnty"
$ws.Range("C3").Value = "Value"
$ws.Range("D3").Value = "Cmp name"
$ws.Range("E3").Value = "Footprint"
$ws.Range("F3").Value = "Vendor"
$ws.Range("G3").Value = "Vendor part number"
$ws.Range("H3").Value = "Received"

# Row 5: C1 (100nF capacitor)
$ws.Range("A5").Value = "C1, "
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "100nF"
$ws.Range("D5").Value = "C_Small"
$ws.Range("E5").Value = "Capacitor_SMD:C_0603_1608Metric"
$ws.Range("H5").Value = "Y"
$ws.Range("I5").Value = "SMT Lab"

# Row 7: C2 (10uF capacitor)
$ws.Range("A7").Value = "C2, "
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "10uF"
$ws.Range("D7").Value = "C_Small"
$ws.Range("E7").Value = "Capacitor_SMD:C_0603_1608Metric"
$ws.Range("H7").Value = "Y"
$ws.Range("I7").Value = "SMT Lab"

# Row 9: H1, H2 (mounting holes)
$ws.Range("A9").Value = "H1, H2, "
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "MountingHole"
$ws.Range("D9").Value = "MountingHole"
$ws.Range("E9").Value = "MountingHole:MountingHole_2.7mm_M2.5_ISO7380"
$ws.Range("F9").Value = "N/A"
$ws.Range("G9").Value = "N/A"
$ws.Range("H9").Value = "N/A"

# Row 11: J1 (connector)
$ws.Range("A11").Value = "J1, "
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "BM04B-SURS-TFLFSN"
$ws.Range("D11").Value = "BM04B-SURS-TFLFSN"
$ws.Range("E11").Value = "BM04B:BM04B-SURS-TFLFSN"
$ws.Range("F11").Value = "LCSC"
$ws.Range("G11").Value = "C495551"
$ws.Range("H11").Value = "Y"

# Row 13: R1, R2 (0R resistors)
$ws.Range("A13").Value = "R1, R2, "
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = "0R"
$ws.Range("D13").Value = "0R"
$ws.Range("E13").Value = "Resistor_SMD:R_0603_1608Metric_Pad0.98x0.95mm_HandSolder"
$ws.Range("F13").Value = "LCSC"
$ws.Range("G13").Value = "C17168"
$ws.Range("H13").Value = "Y"

# Row 15: U1 (AS5600 encoder IC)
$ws.Range("A15").Value = "U1, "
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "AS5600-ASOT"
$ws.Range("D15").Value = "AS5600-ASOT"
$ws.Range("E15").Value = "AS5600:SOIC127P600X175-8N"
$ws.Range("F15").Value = "Digikey"
$ws.Range("G15").Value = "AS5600-ASOTCT-ND"
$ws.Range("H15").Value = "Y"

# Column widths matching the new BOM layout (values chosen so the engine's
# internal character->pixel rounding lands on the closest achievable width
# to the target column widths).
$ws.Columns.Item(1).ColumnWidth = 17.166666666666668
$ws.Columns.Item(2).ColumnWidth = 4.5
$ws.Columns.Item(3).ColumnWidth = 18.666666666666668
$ws.Columns.Item(4).ColumnWidth = 18.666666666666668
$ws.Columns.Item(5).ColumnWidth = 58.666666666666664
$ws.Columns.Item(7).ColumnWidth = 18.5

$ws.Range("I7").Select()
